# Update the "取得日時" (acquisition timestamp) column for all existing
# data rows on the "ランサーズ" sheet from 2025-11-05 01:21:02 to
# 2025-11-05 01:51:31 (append run timestamp refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-11-05 01:21:02"
$newTimestamp = "2025-11-05 01:51:31"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
